$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows whose MarketObjects list changed
$ws.Range("B86").Value = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.QPROBABILITY']"
$ws.Range("B87").Value = "['BTCUSD.SPOT']"

# New rows appended at the bottom of the table
$newRows = @(
    @{ Row = 88; Date = "2025-09-07"; Objs = "['BTCUSD.SPOT']" },
    @{ Row = 89; Date = "2025-09-08"; Objs = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.QPROBABILITY']" },
    @{ Row = 90; Date = "2025-09-09"; Objs = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.QPROBABILITY']" },
    @{ Row = 91; Date = "2025-09-10"; Objs = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.QPROBABILITY']" },
    @{ Row = 92; Date = "2025-09-11"; Objs = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.QPROBABILITY']" },
    @{ Row = 93; Date = "2025-09-12"; Objs = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT']" },
    @{ Row = 94; Date = "2025-09-13"; Objs = "[]" },
    @{ Row = 95; Date = "2025-09-14"; Objs = "[]" }
)

foreach ($r in $newRows) {
    $dateCell = $ws.Cells.Item($r.Row, 1)
    $dateCell.Value = "'" + $r.Date
    $dateCell.Style = "Normal"
    $ws.Cells.Item($r.Row, 2).Value = $r.Objs
}
